# Update cryptocurrency price/volume data (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.379.26"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.13%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.689.60"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.30%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "680.69"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.79%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "159.44"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.43%  "

$ws.Range("E7").Value = "  +0.04%  "

$ws.Range("E8").Value = "  -1.21%  "

$ws.Range("E9").Value = "  -1.44%  "

$ws.Range("E10").Value = "  -3.71%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.438"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.87%  "

$ws.Range("E12").Value = "  -3.31%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.311.09"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.24%  "

$ws.Range("E14").Value = "  -3.19%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.689.74"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.30%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "69.353.00"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.24%  "

$ws.Range("E18").Value = "  -1.91%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.43"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.89%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "468.84"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.79%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.01"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.06%  "

$ws.Range("E22").Value = "  -2.28%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "79.89"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.27%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.835.79"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.21%  "

$ws.Range("E25").Value = "  -0.06%  "

$ws.Range("E26").Value = "  -6.35%  "

$ws.Range("E27").Value = "  -4.32%  "

$ws.Range("E28").Value = "  -4.27%  "

$ws.Range("E29").Value = "  -2.00%  "

$ws.Range("E30").Value = "  -4.18%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.64"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.54%  "

$ws.Range("E32").Value = "  -4.55%  "

$ws.Range("E33").Value = "  +0.00%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "26.95"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.69%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.679.12"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.36%  "

$ws.Range("E36").Value = "  -5.14%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "8.29"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.83%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.22"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.56%  "

$ws.Range("E39").Value = "  +0.00%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.26"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.26%  "

$ws.Range("E41").Value = "  -0.10%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0907"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.12%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "171.78"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.01%  "

$ws.Range("E44").Value = "  -1.15%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "47.61"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.98%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "28.40"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -5.80%  "

$ws.Range("E47").Value = "  -4.36%  "

$ws.Range("B48").Value = "SuiNetwork"
$ws.Range("C48").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.11"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.63%  "

$ws.Range("B49").Value = "FLOKI"
$ws.Range("C49").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.000277"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.48%  "

$ws.Range("E50").Value = "  -5.15%  "

$ws.Range("E51").Value = "  -3.39%  "
